$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.559.18"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.159.97"
$ws.Range("E3").Value = "  +2.95%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.31"
$ws.Range("E5").Value = "  +2.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.09"
$ws.Range("E6").Value = "  +3.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +9.30%  "
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("E10").Value = "  +3.75%  "
$ws.Range("E11").Value = "  +5.50%  "
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.703.30"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.77"
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("E15").Value = "  +7.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.613.60"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.147.80"
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("E18").Value = "  +6.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.05"
$ws.Range("E19").Value = "  +5.31%  "
$ws.Range("E20").Value = "  +6.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.07"
$ws.Range("E21").Value = "  +7.38%  "
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.12"
$ws.Range("E24").Value = "  +2.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.514"
$ws.Range("E25").Value = "  +3.53%  "
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.10"
$ws.Range("E28").Value = "  +14.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0873"
$ws.Range("E29").Value = "  +4.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.18"
$ws.Range("E30").Value = "  +5.99%  "
$ws.Range("E31").Value = "  +2.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.00"
$ws.Range("E32").Value = "  +4.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.20"
$ws.Range("E33").Value = "  +8.75%  "
$ws.Range("E34").Value = "  +5.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.32"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("E36").Value = "  +4.63%  "
$ws.Range("E37").Value = "  +14.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.32"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.67"
$ws.Range("E39").Value = "  +7.47%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.636.48"
$ws.Range("E40").Value = "  +9.34%  "
$ws.Range("E41").Value = "  +4.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.17"
$ws.Range("E42").Value = "  +4.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.74"
$ws.Range("E43").Value = "  +5.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.709"
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("E45").Value = "  +8.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.200.93"
$ws.Range("E47").Value = "  +2.97%  "
$ws.Range("E48").Value = "  +11.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.983"
$ws.Range("E49").Value = "  +5.30%  "
$ws.Range("E50").Value = "  +4.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.27"
$ws.Range("E51").Value = "  +5.20%  "
